$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (sheet1) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 7977.3213
$ws.Range("I33").Value = 10415.096
$ws.Range("K33").Value = 10415.096
$ws.Range("M33").Value = -10186.096
$ws.Range("H95").Value = 277916.34
$ws.Range("J95").Value = 277916.34
$ws.Range("L95").Value = 277916.34
$ws.Range("N95").Value = -283408.34
$ws.Range("H105").Value = 364201.25
$ws.Range("J105").Value = 364201.25
$ws.Range("L105").Value = 364201.25
$ws.Range("N105").Value = -371189.25
$ws.Range("H112").Value = 35717070
$ws.Range("I112").Value = 2373
$ws.Range("K112").Value = 7119
$ws.Range("M112").Value = -6011
$ws.Range("H135").Value = 1116.3334
$ws.Range("I135").Value = 766.3333
$ws.Range("K135").Value = 6896.9997
$ws.Range("M135").Value = -4361.9997
$ws.Range("H137").Value = 18512.28
$ws.Range("J137").Value = 13035.571
$ws.Range("L137").Value = 39106.713
$ws.Range("N137").Value = -44206.713

# ---- Sheet: ARM (sheet2) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3828.6843
$ws.Range("I61").Value = 3605.0667
$ws.Range("K61").Value = 3605.0667
$ws.Range("M61").Value = -3393.0667
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").Value = $null
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 50000
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 50000
$ws.Range("M105").Value = $null
$ws.Range("N105").Value = -56988
$ws.Range("H136").Value = 3828.6843
$ws.Range("I136").Value = 3605.0667
$ws.Range("K136").Value = 10815.2001
$ws.Range("M136").Value = -8265.2001

# ---- Sheet: BSM (sheet3) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 27332.584
$ws.Range("J86").Value = 42641.2
$ws.Range("L86").Value = 42641.2
$ws.Range("N86").Value = -44887.2
$ws.Range("H89").Value = 27332.584
$ws.Range("J89").Value = 42641.2
$ws.Range("L89").Value = 213206
$ws.Range("N89").Value = -224438
$ws.Range("H107").Value = 2936.125
$ws.Range("I107").Value = 2614.2307
$ws.Range("K107").Value = 2614.2307
$ws.Range("M107").Value = -694.2307000000001
$ws.Range("H134").Value = 841.6539
$ws.Range("I134").Value = 796.12
$ws.Range("J134").Value = 1980
$ws.Range("K134").Value = 2388.36
$ws.Range("L134").Value = 5940
$ws.Range("M134").Value = 146.6399999999999
$ws.Range("N134").Value = -11010

# ---- Sheet: CRP (sheet4) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 406894.5
$ws.Range("J43").Value = 406894.5
$ws.Range("L43").Value = 406894.5
$ws.Range("N43").Value = -407262.5
$ws.Range("H101").Value = 406894.5
$ws.Range("J101").Value = 406894.5
$ws.Range("L101").Value = 406894.5
$ws.Range("N101").Value = -413384.5
$ws.Range("H131").Value = 60100
$ws.Range("J131").Value = 60100
$ws.Range("L131").Value = 60100
$ws.Range("N131").Value = -70180
$ws.Range("H134").Value = 2204
$ws.Range("I134").Value = 2152.652
$ws.Range("K134").Value = 6457.956
$ws.Range("M134").Value = -3922.956
$ws.Range("H141").Value = 89621.78999999999
$ws.Range("J141").Value = 89621.78999999999
$ws.Range("L141").Value = 89621.78999999999
$ws.Range("N141").Value = -99981.78999999999

# ---- Sheet: CUL (sheet5) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4000
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").Value = $null
$ws.Range("H83").Value = 4000
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").Value = $null
$ws.Range("H128").Value = 159900
$ws.Range("I128").Value = 159900
$ws.Range("K128").Value = 479700
$ws.Range("M128").Value = -474720
$ws.Range("H137").Value = 1697.1666
$ws.Range("J137").Value = 2031
$ws.Range("L137").Value = 6093
$ws.Range("N137").Value = -16293

# ---- Sheet: GSM (sheet6) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 9500
$ws.Range("J39").Value = 9500
$ws.Range("L39").Value = 9500
$ws.Range("N39").Value = -10564
$ws.Range("H132").Value = 1465
$ws.Range("I132").Value = 1465
$ws.Range("K132").Value = 4395
$ws.Range("M132").Value = -1865

# ---- Sheet: LTW (sheet7) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5270.2
$ws.Range("J7").Value = 7499.75
$ws.Range("L7").Value = 7499.75
$ws.Range("N7").Value = -7723.75
$ws.Range("H40").Value = 27780580
$ws.Range("I40").Value = 41668370
$ws.Range("J40").Value = 5004
$ws.Range("K40").Value = 41668370
$ws.Range("L40").Value = 5004
$ws.Range("M40").Value = -41668234
$ws.Range("N40").Value = -5276
$ws.Range("H100").Value = 3609.8108
$ws.Range("I100").Value = 3308.423
$ws.Range("J100").Value = 4322.1816
$ws.Range("K100").Value = 3308.423
$ws.Range("L100").Value = 4322.1816
$ws.Range("M100").Value = -2767.423
$ws.Range("N100").Value = -5404.1816
$ws.Range("H106").Value = 19632
$ws.Range("J106").Value = 19632
$ws.Range("L106").Value = 19632
$ws.Range("N106").Value = -22156
$ws.Range("H122").Value = 5018.75
$ws.Range("I122").Value = 4309.375
$ws.Range("J122").Value = 6437.5
$ws.Range("K122").Value = 12928.125
$ws.Range("L122").Value = 19312.5
$ws.Range("M122").Value = -10478.125
$ws.Range("N122").Value = -24212.5
$ws.Range("H126").Value = 5270.2
$ws.Range("J126").Value = 7499.75
$ws.Range("L126").Value = 22499.25
$ws.Range("N126").Value = -27439.25
$ws.Range("H131").Value = 89998.5
$ws.Range("J131").Value = 89998.5
$ws.Range("L131").Value = 89998.5
$ws.Range("N131").Value = -100078.5
$ws.Range("H136").Value = 2678.9092
$ws.Range("I136").Value = 2521
$ws.Range("J136").Value = 5995
$ws.Range("K136").Value = 7563
$ws.Range("L136").Value = 17985
$ws.Range("M136").Value = -5013
$ws.Range("N136").Value = -23085

# ---- Sheet: WVR (sheet8) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6333
$ws.Range("I62").Value = 4999.5
$ws.Range("K62").Value = 4999.5
$ws.Range("M62").Value = -4375.5
$ws.Range("H65").Value = 6333
$ws.Range("I65").Value = 4999.5
$ws.Range("K65").Value = 24997.5
$ws.Range("M65").Value = -21877.5
$ws.Range("H95").Value = 44999
$ws.Range("J95").Value = 44998
$ws.Range("L95").Value = 44998
$ws.Range("N95").Value = -50490
$ws.Range("H104").Value = 9499.5
$ws.Range("J104").Value = 9499.5
$ws.Range("L104").Value = 9499.5
$ws.Range("N104").Value = -16487.5
$ws.Range("H105").Value = 43653.5
$ws.Range("J105").Value = 43653.5
$ws.Range("L105").Value = 43653.5
$ws.Range("N105").Value = -50641.5
$ws.Range("H122").Value = 30970.277
$ws.Range("I122").Value = 33747.875
$ws.Range("K122").Value = 101243.625
$ws.Range("M122").Value = -98793.625
$ws.Range("H136").Value = 3749.5476
$ws.Range("I136").Value = 2619.8518
$ws.Range("K136").Value = 7859.555399999999
$ws.Range("M136").Value = -5309.555399999999
